$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135 (shifts rows 135-229 down to 136-230)
$ws.Rows.Item(135).Insert()

# The row that was old row135 is now row136. Copy its values into new row135,
# then overwrite the Fecha (D) and Volumen (J) columns with the new data.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(135, $c).Value = $ws.Cells.Item(136, $c).Value2
}
$ws.Cells.Item(135, 4).Value = Get-Date -Year 2022 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(135, 10).Value = 200

# Copy the date number format from row136 D cell to row135 D cell
$ws.Cells.Item(135, 4).NumberFormat = $ws.Cells.Item(136, 4).NumberFormat
